# The lookup table of Cm_d_alpha values was corrected: every cell in the
# A1:H14 block moves from -1 to -3 (the "w" angular-rate correction bug fix
# referenced in the commit message changed the constant being used here).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# Bulk-update the whole 8x14 table of values from -1 to -3.
$ws.Range("A1:H14").Value = -3

# Reset the view: clear the scrolled-down topLeftCell="A10" state and select
# the full table range A1:H14 (previously A15:H24 was selected).
$ws.Range("A1:H14").Select() | Out-Null
